# Apply updated crypto market data (prices, 1h volume %, and the shifted
# coin rows from LEO onward) per the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.895.14'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.988.38'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.32'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.69'
$ws.Range('E6').Value = '  -2.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +2.45%  '
$ws.Range('D9').Value = '3.002.08'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('E11').Value = '  -3.97%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '3.510.37'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '61.948.70'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.87'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '2.992.59'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.16'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.50'
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.04'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.116.48'
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.469'
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.190'
$ws.Range('E27').Value = '  +1.88%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0924'
$ws.Range('E29').Value = '  -3.99%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.24'
$ws.Range('E30').Value = '  -5.34%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.72'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.47'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '160.37'
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.59'
$ws.Range('E35').Value = '  -2.61%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.91'
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.06'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.27'
$ws.Range('E38').Value = '  -2.10%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.50'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.90'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').Value = '2.415.95'
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.04'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.673'
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0591'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.31'
$ws.Range('E46').Value = '  +6.90%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.996'
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0244'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '269.01'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0949'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.53'
$ws.Range('E51').Value = '  -2.90%  '
